# BUG: Don't extract header names if none specified (#23703)
# Adds a new worksheet "index_col_none" used as a regression fixture for
# pandas' read_excel(index_col=None) with a MultiIndex column header.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab strip (sheetId 13 / rIdN, matching a freshly-minted sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "index_col_none"

# Two header rows: a MultiIndex column header ("A","A","B","B") / ("key","val","key","val")
$newSheet.Range("A1").Value = "A"
$newSheet.Range("B1").Value = "A"
$newSheet.Range("C1").Value = "B"
$newSheet.Range("D1").Value = "B"
$newSheet.Range("A2").Value = "key"
$newSheet.Range("B2").Value = "val"
$newSheet.Range("C2").Value = "key"
$newSheet.Range("D2").Value = "val"

# Two data rows
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = 2
$newSheet.Cells.Item(3, 3).Value = 3
$newSheet.Cells.Item(3, 4).Value = 4
$newSheet.Cells.Item(4, 1).Value = 1
$newSheet.Cells.Item(4, 2).Value = 2
$newSheet.Cells.Item(4, 3).Value = 3
$newSheet.Cells.Item(4, 4).Value = 4

# Center-align the numeric data rows first (creates the new plain/centered
# style), then bold + center the two header rows (reuses the workbook's
# existing bold+centered header style) - matches the order the style table
# ends up in.
$dataRange = $newSheet.Range("A3:D4")
$dataRange.HorizontalAlignment = -4108

$headerRange = $newSheet.Range("A1:D2")
$headerRange.HorizontalAlignment = -4108
$headerRange.Font.Bold = $true
